$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E4").Value = "2016-03-21 22:42:05"
$wsZhCn.Range("H4").Value = "2016-03-21 22:42:28"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E4").Value = "2016-03-21 22:42:09"
$wsDeDe.Range("H4").Value = "2016-03-21 22:42:34"
